$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTasks = @(
    "Make Constructor of Entities for empty collections of relations",
    "Learn about Reqgular Expressions",
    "Make custom Data Annotation attribute for Image Jpg/PNG format",
    "Give validation Error Messages to all Entities",
    "Solve Profile Image Dat Issue",
    "Error Code while registeration if error comes , don’t pass empty userDto",
    "Truncate spaces from starting and ending before saving in database or checkit from frontend while sending it to server",
    "Validation for image jpg/png"
)

$row = 35
foreach ($task in $newTasks) {
    $ws.Cells.Item($row, 2).Value = $task
    $row += 2
}

$ws.Range("E49").Select()
